$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two order lines exceeding the max capacity restriction per station:
#  - Row 11: 4501032020 / P217 / E/S EL ROSARIO / G-PRIX GASOHOL 97 line
#  - Row 43: 4501034585 / P285 / E/S REPÚBLICA / G-PRIX GASOHOL 97 line
# Delete the higher row number first so the second index is not shifted.
$ws.Rows.Item(43).Delete()
$ws.Rows.Item(11).Delete()
